$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Views" (F) column data values for rows 2-6 and 8, while
# preserving the existing percentage number-format/style on those cells.
$ws.Range("F2:F6").ClearContents()
$ws.Range("F8").ClearContents()

# Remove the stray comment text that lived in F9 (and its backing shared
# string) - the cell becomes completely empty.
$ws.Range("F9").ClearContents()
